$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos price/volume figures (Price column D, Volume(1h) column E).
# Values are written as text (matching the source sheet's inline-string cells,
# e.g. "24.931.49", "0.9970") so Excel does not reinterpret them as numbers.
$updates = @(
    @{ Cell = 'D2'; Value = '24.931.49' },
    @{ Cell = 'E2'; Value = '  +0.96%  ' },
    @{ Cell = 'D3'; Value = '1.661.53' },
    @{ Cell = 'E3'; Value = '  -1.99%  ' },
    @{ Cell = 'D4'; Value = '0.9970' },
    @{ Cell = 'E4'; Value = '  -1.09%  ' },
    @{ Cell = 'D5'; Value = '316.95' },
    @{ Cell = 'E5'; Value = '  +1.89%  ' },
    @{ Cell = 'D6'; Value = '0.9972' },
    @{ Cell = 'E6'; Value = '  -0.59%  ' },
    @{ Cell = 'D7'; Value = '0.3637' },
    @{ Cell = 'E7'; Value = '  -2.21%  ' },
    @{ Cell = 'D8'; Value = '46.91' },
    @{ Cell = 'E8'; Value = '  -4.42%  ' },
    @{ Cell = 'D9'; Value = '0.3270' },
    @{ Cell = 'E9'; Value = '  -3.81%  ' },
    @{ Cell = 'D10'; Value = '1.141' },
    @{ Cell = 'E10'; Value = '  -5.00%  ' },
    @{ Cell = 'D11'; Value = '0.07071' },
    @{ Cell = 'E11'; Value = '  -4.78%  ' },
    @{ Cell = 'D12'; Value = '0.9985' },
    @{ Cell = 'E12'; Value = '  -0.60%  ' },
    @{ Cell = 'D13'; Value = '6.055' },
    @{ Cell = 'E13'; Value = '  -3.80%  ' },
    @{ Cell = 'D14'; Value = '19.60' },
    @{ Cell = 'E14'; Value = '  -5.88%  ' },
    @{ Cell = 'D15'; Value = '1.663.14' },
    @{ Cell = 'E15'; Value = '  -1.94%  ' },
    @{ Cell = 'D16'; Value = '6.626' },
    @{ Cell = 'E16'; Value = '  -4.68%  ' },
    @{ Cell = 'D17'; Value = '0.00001050' },
    @{ Cell = 'E17'; Value = '  -5.80%  ' },
    @{ Cell = 'D18'; Value = '0.06629' },
    @{ Cell = 'E18'; Value = '  -1.00%  ' },
    @{ Cell = 'D19'; Value = '0.9971' },
    @{ Cell = 'E19'; Value = '  -0.52%  ' },
    @{ Cell = 'D20'; Value = '79.35' },
    @{ Cell = 'E20'; Value = '  -4.18%  ' },
    @{ Cell = 'D21'; Value = '5.918' },
    @{ Cell = 'E21'; Value = '  -5.98%  ' },
    @{ Cell = 'D22'; Value = '15.78' },
    @{ Cell = 'E22'; Value = '  -7.69%  ' },
    @{ Cell = 'E23'; Value = '  -2.50%  ' },
    @{ Cell = 'D24'; Value = '24.798.65' },
    @{ Cell = 'E24'; Value = '  +0.47%  ' },
    @{ Cell = 'D25'; Value = '2.430' },
    @{ Cell = 'E25'; Value = '  -0.70%  ' },
    @{ Cell = 'D26'; Value = '2.394' },
    @{ Cell = 'E26'; Value = '  -12.96%  ' },
    @{ Cell = 'D27'; Value = '148.54' },
    @{ Cell = 'E27'; Value = '  +0.15%  ' },
    @{ Cell = 'E28'; Value = '  -7.44%  ' },
    @{ Cell = 'D29'; Value = '1.224' },
    @{ Cell = 'E29'; Value = '  +0.40%  ' },
    @{ Cell = 'D30'; Value = '1.847.93' },
    @{ Cell = 'E30'; Value = '  -1.90%  ' },
    @{ Cell = 'D31'; Value = '125.99' },
    @{ Cell = 'E31'; Value = '  -3.88%  ' },
    @{ Cell = 'D32'; Value = '4.086' },
    @{ Cell = 'E32'; Value = '  -2.97%  ' },
    @{ Cell = 'D33'; Value = '5.833' },
    @{ Cell = 'E33'; Value = '  -12.55%  ' },
    @{ Cell = 'D34'; Value = '0.08446' },
    @{ Cell = 'E34'; Value = '  -2.59%  ' },
    @{ Cell = 'D35'; Value = '1.679' },
    @{ Cell = 'E35'; Value = '  -3.76%  ' },
    @{ Cell = 'D36'; Value = '12.30' },
    @{ Cell = 'E36'; Value = '  -8.93%  ' },
    @{ Cell = 'D37'; Value = '1.279' },
    @{ Cell = 'E37'; Value = '  +1.91%  ' },
    @{ Cell = 'D38'; Value = '5.209' },
    @{ Cell = 'E38'; Value = '  -5.29%  ' },
    @{ Cell = 'D39'; Value = '0.06048' },
    @{ Cell = 'E39'; Value = '  -8.04%  ' },
    @{ Cell = 'D40'; Value = '0.02238' },
    @{ Cell = 'E40'; Value = '  -6.31%  ' },
    @{ Cell = 'D41'; Value = '0.2069' },
    @{ Cell = 'E41'; Value = '  -6.15%  ' },
    @{ Cell = 'D42'; Value = '8.234' },
    @{ Cell = 'E42'; Value = '  -8.58%  ' },
    @{ Cell = 'D43'; Value = '0.9970' },
    @{ Cell = 'E43'; Value = '  -0.60%  ' },
    @{ Cell = 'D44'; Value = '0.5920' },
    @{ Cell = 'E44'; Value = '  -6.90%  ' },
    @{ Cell = 'D45'; Value = '3.815' },
    @{ Cell = 'E45'; Value = '  +0.08%  ' },
    @{ Cell = 'D46'; Value = '12.73' },
    @{ Cell = 'E46'; Value = '  -6.71%  ' },
    @{ Cell = 'D47'; Value = '0.5646' },
    @{ Cell = 'E47'; Value = '  -6.66%  ' },
    @{ Cell = 'D48'; Value = '125.11' },
    @{ Cell = 'E48'; Value = '  -2.61%  ' },
    @{ Cell = 'D49'; Value = '1.952' },
    @{ Cell = 'E49'; Value = '  -7.08%  ' },
    @{ Cell = 'D50'; Value = '0.07014' },
    @{ Cell = 'E50'; Value = '  -2.96%  ' },
    @{ Cell = 'D51'; Value = '1.196' },
    @{ Cell = 'E51'; Value = '  -1.12%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
